$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value2 = '[{"factor": ["thioredoxin concentration", "serum"], "outcome": ["in-hospital major adverse event"]}, {"factor": ["glasgow coma scale score", "admission"], "outcome": ["in-hospital major adverse event"]}]'
$ws.Range("G4").Value2 = '[{"factor": ["mif", "macrophage migration inhibitory factor"], "outcome": ["severity", "inflammation", "clinical outcome", "trauma"]}]'
$ws.Range("G5").Value2 = '[{"factor": ["glasgow coma scale score"], "outcome": ["glasgow outcome scale", "gos"]}, {"factor": ["artificial ventilation"], "outcome": ["mortality rate"]}, {"factor": ["monitoring", "intracranial"], "outcome": ["outcome", "long-term", "good"]}]'
$ws.Range("G6").Value2 = '[{"factor": ["impact", "calculator", "prognosis"], "outcome": ["patient", "elderly", "severe", "outcome", "tbi"]}, {"factor": ["expect risk", "fatal", "outcome"], "outcome": ["risk", "conservative treatment"]}, {"factor": ["outcome", "unfavorable", "predict risk"], "outcome": ["actual", "rate", "risk", "outcome", "unfavorable", "predict"]}]'
$ws.Range("G8").Value2 = '[{"factor": ["compute tomography", "ct"], "outcome": ["death", "patient", "tbi"]}, {"factor": ["marshall and rotterdam scoring system"], "outcome": ["death", "patient", "tbi"]}, {"factor": ["basal cistern absence"], "outcome": ["death", "patient", "tbi"]}, {"factor": ["positive midline shift"], "outcome": ["death", "patient", "tbi"]}, {"factor": ["hemorrhagic mass volume"], "outcome": ["death", "patient", "tbi"]}, {"factor": ["intraventricular", "subarachnoid hemorrhage"], "outcome": ["death", "patient", "tbi"]}]'
$ws.Range("G9").Value2 = '[{"factor": ["apache ii"], "outcome": ["icu-treated", "mortality", "patient", "tbi", "six-month"]}, {"factor": ["sap ii"], "outcome": ["icu-treated", "mortality", "patient", "tbi", "six-month"]}, {"factor": ["sofa"], "outcome": ["icu-treated", "mortality", "patient", "tbi", "six-month"]}, {"factor": ["age"], "outcome": ["icu-treated", "mortality", "patient", "tbi", "six-month"]}, {"factor": ["glasgow coma scale"], "outcome": ["icu-treated", "mortality", "patient", "tbi", "six-month"]}]'
$ws.Range("G10").Value2 = '[{"factor": ["v/c ratio"], "outcome": ["drs", "gos score", "lcf score"]}]'
$ws.Range("G11").Value2 = '[{"factor": ["timp-1 level", "serum"], "outcome": ["tbi", "patient", "severe", "mortality"]}]'
$ws.Range("G12").Value2 = '[{"factor": ["evidence", "contusion", "mri"], "outcome": ["glasgow outcome scale-extended", "gos-e"]}, {"factor": ["roi", "severely", "fa", "reduce"], "outcome": ["gos-e"]}, {"factor": ["neuropsychiatric history"], "outcome": ["gos-e"]}, {"factor": ["age"], "outcome": ["gos-e"]}, {"factor": ["year of"], "outcome": ["gos-e"]}, {"factor": ["roi", "severely", "fa", "reduce"], "outcome": ["gos-e"]}, {"factor": ["neuropsychiatric history"], "outcome": ["gos-e"]}, {"factor": ["year of"], "outcome": ["gos-e"]}]'
$ws.Range("G13").Value2 = '[{"factor": ["preinjury", "depressive symptom"], "outcome": ["affective/behavioral", "physical problem", "cognitive", "health-related quality-of-life"]}]'
$ws.Range("G15").Value2 = '[{"factor": ["age"], "outcome": ["in-hospital mortality rate"]}, {"factor": ["sex"], "outcome": ["in-hospital mortality rate"]}, {"factor": ["score", "deyo-charlson comorbidity index"], "outcome": ["in-hospital mortality rate"]}, {"factor": ["hospital volume"], "outcome": ["in-hospital mortality rate"]}, {"factor": ["surgeon", "volume"], "outcome": ["in-hospital mortality rate"]}]'
$ws.Range("G16").Value2 = '[{"factor": ["age"], "outcome": ["risk", "icu"]}, {"factor": ["gender"], "outcome": ["risk", "icu"]}, {"factor": ["etiology", "injury"], "outcome": ["risk", "icu"]}]'
$ws.Range("G17").Value2 = '[{"factor": ["concentration", "cerebrospinal fluid", "csf", "amyloid-beta1-42 (abeta42)"], "outcome": ["mortality"]}, {"factor": ["concentration", "plasma", "amyloid-beta1-42 (abeta42)"], "outcome": ["mortality"]}, {"factor": ["abeta42", "change", "concentration", "csf"], "outcome": ["neurological status"]}]'
$ws.Range("G18").Value2 = '[{"factor": ["plasminogen activator receptor", "supar", "soluble", "urokinase"], "outcome": ["traumatic brain injury", "tbi"]}, {"factor": ["glasgow coma scale", "score"], "outcome": ["severity", "tbi"]}, {"factor": ["d-dimer"], "outcome": ["tbi", "patient", "mortality"]}]'
$ws.Range("G19").Value2 = '[{"factor": ["motor score", "gcs"], "outcome": ["month", "mortality"]}, {"factor": ["age"], "outcome": ["month", "mortality"]}, {"factor": ["sex"], "outcome": ["month", "mortality"]}, {"factor": ["injury", "mechanism"], "outcome": ["month", "mortality"]}, {"factor": ["glasgow coma scale"], "outcome": ["month", "mortality"]}, {"factor": ["intubation"], "outcome": ["month", "mortality"]}, {"factor": ["pupil"], "outcome": ["month", "mortality"]}, {"factor": ["systolic blood pressure"], "outcome": ["month", "mortality"]}, {"factor": ["respiratory rate"], "outcome": ["month", "mortality"]}, {"factor": ["body temperature"], "outcome": ["month", "mortality"]}, {"factor": ["arterial", "ph"], "outcome": ["month", "mortality"]}, {"factor": ["arterial partial pressure", "carbon dioxide"], "outcome": ["month", "mortality"]}, {"factor": ["arterial partial pressure"], "outcome": ["month", "mortality"]}, {"factor": ["serum sodium"], "outcome": ["month", "mortality"]}, {"factor": ["serum potassium"], "outcome": ["month", "mortality"]}, {"factor": ["serum chloride"], "outcome": ["month", "mortality"]}, {"factor": ["serum calcium"], "outcome": ["month", "mortality"]}, {"factor": ["serum glucose"], "outcome": ["month", "mortality"]}, {"factor": ["urea nitrogen"], "outcome": ["month", "mortality"]}, {"factor": ["creatinine"], "outcome": ["month", "mortality"]}, {"factor": ["ratio", "international"], "outcome": ["month", "mortality"]}]'
$ws.Range("G22").Value2 = '[{"factor": ["systolic blood pressure", "sbp"], "outcome": ["mortality"]}]'
$ws.Range("G24").Value2 = '[{"factor": ["tsp-1"], "outcome": ["outcome", "unfavorable", "1-week", "mortality"]}]'
$ws.Range("G25").Value2 = '[{"factor": ["brain-derived neurotrophic factor (bdnf)", "plasma level"], "outcome": ["icu", "mortality", "patient", "severe", "intensive care unit", "tbi"]}]'
$ws.Range("G26").Value2 = '[{"factor": ["crash-ct model"], "outcome": ["death", "days"]}, {"factor": ["age"], "outcome": ["patient", "performance", "model", "older"]}, {"factor": ["glasgow coma scale score"], "outcome": ["discrimination", "model"]}, {"factor": ["hosmer-lemeshow p value"], "outcome": ["calibration", "model"]}]'
$ws.Range("G27").Value2 = '[{"factor": ["time to death"], "outcome": ["withdrawal", "life-sustaining", "therapy"]}, {"factor": ["glasgow coma scale", "score"], "outcome": ["mortality"]}, {"factor": ["score", "ais", "head abbreviate injury scale"], "outcome": ["mortality"]}, {"factor": ["comorbiditie", "multiple"], "outcome": ["mortality"]}, {"factor": ["subarachnoid hemorrhage", "traumatic"], "outcome": ["mortality"]}, {"factor": ["intracerebral mass lesion"], "outcome": ["mortality"]}, {"factor": ["brainstem lesion"], "outcome": ["mortality"]}, {"factor": ["absent", "sign of compress", "basal cistern"], "outcome": ["mortality"]}]'
$ws.Range("G28").Value2 = '[{"factor": ["core model"], "outcome": ["mortality"]}, {"factor": ["extend model"], "outcome": ["mortality"]}, {"factor": ["lab model"], "outcome": ["mortality"]}, {"factor": ["core model"], "outcome": ["outcome", "unfavorable"]}, {"factor": ["extend model"], "outcome": ["outcome", "unfavorable"]}, {"factor": ["lab model"], "outcome": ["outcome", "unfavorable"]}]'
$ws.Range("G29").Value2 = '[{"factor": ["glasgow coma scale", "score"], "outcome": ["outcome", "unfavorable", "1-week", "mortality"]}]'
$ws.Range("G30").Value2 = '[{"factor": ["gcs", "glasgow coma scale"], "outcome": ["severity", "tbi"]}, {"factor": ["pta", "duration", "post-traumatic amnesia"], "outcome": ["olfactory problem", "tbi"]}]'
$ws.Range("G31").Value2 = '[{"factor": ["level", "il-6"], "outcome": ["septic", "development"]}, {"factor": ["c-reactive protein level"], "outcome": ["multiple organ dysfunction", "development"]}]'
$ws.Range("G32").Value2 = '[{"factor": ["rotterdam"], "outcome": ["death", "hospital discharge", "weeks"]}, {"factor": ["age"], "outcome": ["death", "hospital discharge", "weeks"]}, {"factor": ["sex"], "outcome": ["death", "hospital discharge", "weeks"]}, {"factor": ["glasgow coma scale score"], "outcome": ["death", "hospital discharge", "weeks"]}]'
$ws.Range("G33").Value2 = '[{"factor": ["glucose", "serum", "admission", "level of"], "outcome": ["outcome", "traumatic brain injury", "patient", "severe"]}]'
$ws.Range("G36").Value2 = '[{"factor": ["il-6"], "outcome": ["year", "favorable", "gos"]}, {"factor": ["pg"], "outcome": ["year", "favorable", "gos"]}, {"factor": ["gfap"], "outcome": ["unfavorable", "year", "score", "gos"]}, {"factor": ["pg"], "outcome": ["survival status", "year"]}, {"factor": ["gfap"], "outcome": ["survival status", "year"]}, {"factor": ["il-6"], "outcome": ["survival status", "year"]}]'
$ws.Range("G38").Value2 = '[{"factor": ["glasgow coma scale"], "outcome": ["mortality"]}, {"factor": ["mechanical ventilation"], "outcome": ["neurological"]}, {"factor": ["blood transfusion"], "outcome": ["neurological"]}, {"factor": ["neurosurgical intervention"], "outcome": ["neurological"]}, {"factor": ["concomitant", "injury"], "outcome": ["complication", "non-neurological"]}]'
$ws.Range("G39").Value2 = '[{"factor": ["traumatic brain injury"], "outcome": ["mortality"]}, {"factor": ["hospital", "systolic blood pressure", "refer"], "outcome": ["mortality"]}]'
